$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Estudos")

# Row 62 corresponds to 22/11/2022 (end of work day).
# Fill in HORA F (C62) and DESCANSO (E62); DIF (D62) and UTEIS (F62) are
# calculated columns in the table and will recompute automatically.
$ws.Range("C62").Value = 0.80555555555555547
$ws.Range("E62").Value = 0.12361111111111112

# Update ASSUNTO (G62) and PRODUÇÃO (H62) text for the day.
$ws.Range("G62").Value = "ESTÁGIO + HARD"
$ws.Range("H62").Value = "Estágio + HARD"

$wb.Application.Calculate()
